$d = $word.ActiveDocument
$lineBreak = [char]11

# 1. Replace "LOQ4031 -  Química Geral I  (Requisito)" with
#    "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)"
$rng1 = $d.Content
$rng1.Find.Execute("LOQ4031 -  Química Geral I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)", 2)

# 2. Remove the whole "LOQ4073 -  Química Geral II  (Requisito)" line (text + its line break)
$rng2 = $d.Content
$rng2.Find.Execute("LOQ4073 -  Química Geral II  (Requisito)" + $lineBreak, $true, $false, $false, $false, $false, $true, 1, $false, "", 1)

# 3. After the "LOQ4095 -  Química Geral Experimental  (Requisito)" line, insert two new
#    requirement lines as separate runs, each with its own trailing line break.
$rng3 = $d.Content
$rng3.Find.Execute("LOQ4095 -  Química Geral Experimental  (Requisito)" + $lineBreak, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ins1 = $rng3.Duplicate
$ins1.Collapse(0)
$ins1.InsertAfter("LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)" + $lineBreak)
$ins1.Collapse(0)

$rng4 = $d.Content
$rng4.Find.Execute("LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)" + $lineBreak, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ins2 = $rng4.Duplicate
$ins2.Collapse(0)
$ins2.InsertAfter("LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)" + $lineBreak)
